# Calendar.xlsx: add a new "GS-QC-6301" class block spanning rows 6-12
# (columns B and F), shrinking the content that used to start at B10/F10
# down to an empty, still-merged placeholder.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$nbsp = [char]0x00A0
$newText = "GS-QC-6301" + $nbsp + "`n" + "09:00-10:30" + $nbsp + "`n"

# The text that used to live at the top of the B10:B14 / F10:F14 merged
# blocks moves up into the new B6:B12 / F6:F12 blocks, so B10/F10 become
# blank (but keep their centered/wrapped style).
$ws.Range("B10").Value = ""
$ws.Range("F10").Value = ""

# Give the new block the same centered + wrap-text styling used elsewhere
# in the grid before merging it.
$ws.Range("B6:B12").Style = $ws.Range("B11").Style
$ws.Range("F6:F12").Style = $ws.Range("F11").Style

$ws.Range("B6:B12").Merge()
$ws.Range("F6:F12").Merge()

$ws.Range("B6").Value = $newText
$ws.Range("F6").Value = $newText

# Merging a wrapped multi-line cell auto-expands the row height; put row 6
# back to the sheet's normal auto height.
$ws.Rows(6).AutoFit()
